# Refresh the cryptocurrency price / 1h-volume table on Sheet1 with the
# latest scrape, including the two row pairs whose ranking order swapped
# (RenderToken <-> FirstDigitalUSD, Stacks <-> EthereumClassic).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.093.46'
$ws.Range("E2").Value = '  +2.13%  '

$ws.Range("D3").Value = '2.416.42'
$ws.Range("E3").Value = '  +2.83%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.37'
$ws.Range("E5").Value = '  +2.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.10'
$ws.Range("E6").Value = '  +4.60%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  +2.29%  '

$ws.Range("D9").Value = '2.414.10'
$ws.Range("E9").Value = '  +2.74%  '

$ws.Range("E10").Value = '  +3.92%  '

$ws.Range("E11").Value = '  -0.65%  '

$ws.Range("E12").Value = '  +1.25%  '

$ws.Range("E13").Value = '  +1.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.22'
$ws.Range("E14").Value = '  +6.28%  '

$ws.Range("E15").Value = '  +8.19%  '

$ws.Range("D16").Value = '2.856.24'
$ws.Range("E16").Value = '  +2.92%  '

$ws.Range("D17").Value = '62.015.06'
$ws.Range("E17").Value = '  +2.19%  '

$ws.Range("D18").Value = '2.417.15'
$ws.Range("E18").Value = '  +2.78%  '

$ws.Range("E19").Value = '  +3.87%  '

$ws.Range("E20").Value = '  +1.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.27'
$ws.Range("E21").Value = '  +1.23%  '

$ws.Range("E22").Value = '  +2.31%  '

$ws.Range("E23").Value = '  +0.19%  '

$ws.Range("E24").Value = '  +5.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.93'
$ws.Range("E25").Value = '  +2.36%  '

$ws.Range("E26").Value = '  +6.98%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '584.01'
$ws.Range("E27").Value = '  +17.58%  '

$ws.Range("D28").Value = '2.537.03'
$ws.Range("E28").Value = '  +2.87%  '

$ws.Range("E29").Value = '  +0.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.36'
$ws.Range("E30").Value = '  +4.28%  '

$ws.Range("D31").Value = '0.0₃0936'
$ws.Range("E31").Value = '  +8.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.47'
$ws.Range("E32").Value = '  +6.78%  '

$ws.Range("E34").Value = '  +3.47%  '

$ws.Range("E35").Value = '  +2.93%  '

$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.68'
$ws.Range("E36").Value = '  +8.39%  '

$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.82'
$ws.Range("E38").Value = '  +3.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.384'
$ws.Range("E39").Value = '  +2.15%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.88'
$ws.Range("E40").Value = '  +3.05%  '

$ws.Range("B41").Value = 'EthereumClassic'
$ws.Range("C41").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.77'
$ws.Range("E41").Value = '  +1.33%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '148.69'
$ws.Range("E42").Value = '  +2.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.76'
$ws.Range("E44").Value = '  +2.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '151.23'
$ws.Range("E45").Value = '  +5.93%  '

$ws.Range("E46").Value = '  +12.42%  '

$ws.Range("E47").Value = '  +2.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0546'
$ws.Range("E48").Value = '  +5.77%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.35'
$ws.Range("E49").Value = '  +6.60%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.588'
$ws.Range("E50").Value = '  +3.45%  '

$ws.Range("E51").Value = '  +1.63%  '
